$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.520.91"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "1.976.98"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'240.85"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'56.24"
$ws.Range("E8").Value = "  +6.78%  "
$ws.Range("D9").Value = "'59.69"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "'0.358"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.0726"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("E12").Value = "  -5.14%  "
$ws.Range("D13").Value = "'0.894"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'14.18"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").Value = "2.265.49"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "1.970.03"
$ws.Range("E17").Value = "  -4.81%  "
$ws.Range("D18").Value = "'17.12"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("D19").Value = "35.343.13"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").Value = "'69.96"
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").Value = "'231.80"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  -4.74%  "
$ws.Range("E26").Value = "  +6.41%  "
$ws.Range("D27").Value = "'163.01"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'9.05"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "'19.37"
$ws.Range("E29").Value = "  -4.84%  "
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  -5.54%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E35").Value = "  -7.28%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "'4.85"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.18"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("E43").Value = "  -5.32%  "
$ws.Range("D44").Value = "'0.0887"
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("D45").Value = "'90.30"
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("D46").Value = "1.370.86"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "'15.33"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").Value = "'45.39"
$ws.Range("E51").Value = "  +1.86%  "
